# Apply updates to the "Inscricoes" sheet of the workbook.
# These changes reflect updated registration counts (Inscritos/Pagos/
# Inscrições homologadas) for a few course rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 3: Pagos (F3) and Inscrições homologadas (H3): 8 -> 9
$ws.Range("F3").Value = 9
$ws.Range("H3").Value = 9

# Row 7: Pagos (F7) and Inscrições homologadas (H7): 12 -> 13
$ws.Range("F7").Value = 13
$ws.Range("H7").Value = 13

# Row 11: Inscritos (E11): 11 -> 12
$ws.Range("E11").Value = 12

# Row 15: Inscritos (E15): 83 -> 84
$ws.Range("E15").Value = 84

# Row 16: Inscritos (E16): 286 -> 290; Pagos (F16) and Inscrições
# homologadas (H16): 80 -> 82
$ws.Range("E16").Value = 290
$ws.Range("F16").Value = 82
$ws.Range("H16").Value = 82
